# Applies numeric cell updates per the target diff (Betfair Back/Lay odds sheet).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 6).Value = 4
$ws.Cells.Item(2, 7).Value = 4.5
$ws.Cells.Item(2, 8).Value = 2.06
$ws.Cells.Item(2, 9).Value = 2.1
$ws.Cells.Item(2, 10).Value = 3.45
$ws.Cells.Item(2, 12).Value = 1.45
$ws.Cells.Item(2, 16).Value = 1.87
$ws.Cells.Item(2, 17).Value = 2.1
$ws.Cells.Item(2, 18).Value = 1.31
$ws.Cells.Item(2, 20).Value = 1.87
$ws.Cells.Item(2, 22).Value = 1.9
$ws.Cells.Item(2, 23).Value = 1.3
# Row 3
$ws.Cells.Item(3, 6).Value = 1.47
$ws.Cells.Item(3, 7).Value = 1.58
$ws.Cells.Item(3, 8).Value = 7.4
$ws.Cells.Item(3, 9).Value = 8.800000000000001
$ws.Cells.Item(3, 10).Value = 4.3
$ws.Cells.Item(3, 11).Value = 5.1
$ws.Cells.Item(3, 20).Value = 1.89
$ws.Cells.Item(3, 22).Value = 1.12
$ws.Cells.Item(3, 23).Value = 2.72
# Row 4
$ws.Cells.Item(4, 7).Value = 2
$ws.Cells.Item(4, 12).Value = 1.39
$ws.Cells.Item(4, 15).Value = 1.32
$ws.Cells.Item(4, 17).Value = 1.95
$ws.Cells.Item(4, 18).Value = 1.39
# Row 5
$ws.Cells.Item(5, 7).Value = 1.36
$ws.Cells.Item(5, 8).Value = 12.5
$ws.Cells.Item(5, 9).Value = 13.5
$ws.Cells.Item(5, 11).Value = 5.6
$ws.Cells.Item(5, 22).Value = 1.08
$ws.Cells.Item(5, 23).Value = 3.8
$ws.Cells.Item(5, 25).Value = 32
$ws.Cells.Item(5, 30).Value = 48
$ws.Cells.Item(5, 31).Value = 280
$ws.Cells.Item(5, 39).Value = 290
# Row 6
$ws.Cells.Item(6, 6).Value = 2.8
$ws.Cells.Item(6, 7).Value = 3.05
$ws.Cells.Item(6, 9).Value = 2.72
$ws.Cells.Item(6, 10).Value = 3.4
$ws.Cells.Item(6, 14).Value = 3.85
$ws.Cells.Item(6, 16).Value = 1.99
$ws.Cells.Item(6, 18).Value = 1.39
$ws.Cells.Item(6, 19).Value = 3.05
$ws.Cells.Item(6, 20).Value = 1.67
$ws.Cells.Item(6, 21).Value = 2.22
$ws.Cells.Item(6, 22).Value = 1.58
$ws.Cells.Item(6, 23).Value = 1.49
$ws.Cells.Item(6, 24).Value = 19.5
$ws.Cells.Item(6, 25).Value = 14.5
$ws.Cells.Item(6, 26).Value = 22
$ws.Cells.Item(6, 27).Value = 46
$ws.Cells.Item(6, 28).Value = 15.5
$ws.Cells.Item(6, 29).Value = 10
$ws.Cells.Item(6, 30).Value = 15
$ws.Cells.Item(6, 31).Value = 34
$ws.Cells.Item(6, 32).Value = 25
$ws.Cells.Item(6, 33).Value = 15.5
$ws.Cells.Item(6, 34).Value = 20
$ws.Cells.Item(6, 35).Value = 46
$ws.Cells.Item(6, 36).Value = 55
$ws.Cells.Item(6, 37).Value = 38
$ws.Cells.Item(6, 38).Value = 48
$ws.Cells.Item(6, 39).Value = 95
$ws.Cells.Item(6, 40).Value = 30
$ws.Cells.Item(6, 41).Value = 26
# Row 7
$ws.Cells.Item(7, 6).Value = 1.58
$ws.Cells.Item(7, 7).Value = 1.69
$ws.Cells.Item(7, 8).Value = 5.6
$ws.Cells.Item(7, 9).Value = 6.8
$ws.Cells.Item(7, 10).Value = 4.2
$ws.Cells.Item(7, 12).Value = 1.29
$ws.Cells.Item(7, 14).Value = 4.3
$ws.Cells.Item(7, 15).Value = 1.24
$ws.Cells.Item(7, 17).Value = 1.71
$ws.Cells.Item(7, 19).Value = 2.8
$ws.Cells.Item(7, 20).Value = 1.8
$ws.Cells.Item(7, 21).Value = 2.02
$ws.Cells.Item(7, 22).Value = 1.17
$ws.Cells.Item(7, 23).Value = 2.44
$ws.Cells.Item(7, 25).Value = 25
$ws.Cells.Item(7, 28).Value = 10.5
$ws.Cells.Item(7, 29).Value = 11
$ws.Cells.Item(7, 30).Value = 1000
$ws.Cells.Item(7, 32).Value = 11
$ws.Cells.Item(7, 33).Value = 11
$ws.Cells.Item(7, 34).Value = 1000
$ws.Cells.Item(7, 36).Value = 16
$ws.Cells.Item(7, 37).Value = 17.5
$ws.Cells.Item(7, 38).Value = 1000
$ws.Cells.Item(7, 40).Value = 8.4
# Row 8
$ws.Cells.Item(8, 12).Value = 1.24
$ws.Cells.Item(8, 21).Value = 2.14
# Row 9
$ws.Cells.Item(9, 6).Value = 2.1
$ws.Cells.Item(9, 7).Value = 2.26
$ws.Cells.Item(9, 8).Value = 3.35
$ws.Cells.Item(9, 9).Value = 3.75
$ws.Cells.Item(9, 10).Value = 3.75
$ws.Cells.Item(9, 11).Value = 4.2
$ws.Cells.Item(9, 12).Value = 1.24
$ws.Cells.Item(9, 14).Value = 5.2
$ws.Cells.Item(9, 15).Value = 1.2
$ws.Cells.Item(9, 17).Value = 1.6
$ws.Cells.Item(9, 18).Value = 1.58
$ws.Cells.Item(9, 19).Value = 2.48
$ws.Cells.Item(9, 20).Value = 1.56
$ws.Cells.Item(9, 21).Value = 2.56
$ws.Cells.Item(9, 22).Value = 1.37
$ws.Cells.Item(9, 23).Value = 1.81
$ws.Cells.Item(9, 24).Value = 27
$ws.Cells.Item(9, 25).Value = 970
$ws.Cells.Item(9, 26).Value = 34
$ws.Cells.Item(9, 27).Value = 70
$ws.Cells.Item(9, 28).Value = 17
$ws.Cells.Item(9, 29).Value = 11.5
$ws.Cells.Item(9, 30).Value = 18.5
$ws.Cells.Item(9, 31).Value = 42
$ws.Cells.Item(9, 32).Value = 20
$ws.Cells.Item(9, 33).Value = 14
$ws.Cells.Item(9, 34).Value = 17.5
$ws.Cells.Item(9, 35).Value = 970
$ws.Cells.Item(9, 36).Value = 34
$ws.Cells.Item(9, 37).Value = 25
$ws.Cells.Item(9, 38).Value = 34
$ws.Cells.Item(9, 39).Value = 70
$ws.Cells.Item(9, 40).Value = 13
$ws.Cells.Item(9, 41).Value = 970
# Row 10
$ws.Cells.Item(10, 6).Value = 1.8
$ws.Cells.Item(10, 7).Value = 1.81
$ws.Cells.Item(10, 10).Value = 4.1
$ws.Cells.Item(10, 11).Value = 4.2
$ws.Cells.Item(10, 12).Value = 1.29
$ws.Cells.Item(10, 17).Value = 1.63
$ws.Cells.Item(10, 18).Value = 1.6
$ws.Cells.Item(10, 19).Value = 2.6
$ws.Cells.Item(10, 21).Value = 2.44
$ws.Cells.Item(10, 23).Value = 2.22
$ws.Cells.Item(10, 40).Value = 8.199999999999999
# Row 11
$ws.Cells.Item(11, 12).Value = 1.3
$ws.Cells.Item(11, 14).Value = 5.4
$ws.Cells.Item(11, 16).Value = 2.5
$ws.Cells.Item(11, 18).Value = 1.61
$ws.Cells.Item(11, 27).Value = 75
$ws.Cells.Item(11, 28).Value = 13
# Row 12
$ws.Cells.Item(12, 12).Value = 1.17
$ws.Cells.Item(12, 16).Value = 4.2
$ws.Cells.Item(12, 18).Value = 2.32
$ws.Cells.Item(12, 19).Value = 1.71
$ws.Cells.Item(12, 31).Value = 12.5
# Row 13
$ws.Cells.Item(13, 8).Value = 27
$ws.Cells.Item(13, 9).Value = 29
$ws.Cells.Item(13, 23).Value = 7.2
$ws.Cells.Item(13, 33).Value = 14.5
$ws.Cells.Item(13, 35).Value = 510
$ws.Cells.Item(13, 37).Value = 16
$ws.Cells.Item(13, 39).Value = 550
$ws.Cells.Item(13, 40).Value = 3.35
# Row 14
$ws.Cells.Item(14, 9).Value = 10.5
$ws.Cells.Item(14, 12).Value = 1.39
$ws.Cells.Item(14, 20).Value = 2.1
$ws.Cells.Item(14, 22).Value = 1.1
$ws.Cells.Item(14, 25).Value = 29
$ws.Cells.Item(14, 27).Value = 410
$ws.Cells.Item(14, 28).Value = 7.8
$ws.Cells.Item(14, 30).Value = 970
$ws.Cells.Item(14, 34).Value = 970
# Row 15
$ws.Cells.Item(15, 7).Value = 7.8
$ws.Cells.Item(15, 12).Value = 1.3
$ws.Cells.Item(15, 17).Value = 1.62
$ws.Cells.Item(15, 19).Value = 2.54
$ws.Cells.Item(15, 20).Value = 1.8
$ws.Cells.Item(15, 21).Value = 2.1
$ws.Cells.Item(15, 23).Value = 1.15
$ws.Cells.Item(15, 25).Value = 10.5
$ws.Cells.Item(15, 32).Value = 65
$ws.Cells.Item(15, 36).Value = 280
$ws.Cells.Item(15, 37).Value = 95
$ws.Cells.Item(15, 38).Value = 85
# Row 16
$ws.Cells.Item(16, 13).Value = 1.03
$ws.Cells.Item(16, 17).Value = 1.75
$ws.Cells.Item(16, 19).Value = 1.75
# Row 17
$ws.Cells.Item(17, 6).Value = 2.68
$ws.Cells.Item(17, 8).Value = 2.68
$ws.Cells.Item(17, 11).Value = 3.85
$ws.Cells.Item(17, 13).Value = 1.06
$ws.Cells.Item(17, 18).Value = 1.36
$ws.Cells.Item(17, 19).Value = 3.25
$ws.Cells.Item(17, 20).Value = 1.69
$ws.Cells.Item(17, 21).Value = 2.18
$ws.Cells.Item(17, 24).Value = 970
$ws.Cells.Item(17, 25).Value = 12.5
$ws.Cells.Item(17, 26).Value = 20
$ws.Cells.Item(17, 27).Value = 44
$ws.Cells.Item(17, 28).Value = 12.5
$ws.Cells.Item(17, 29).Value = 8.4
$ws.Cells.Item(17, 30).Value = 13
$ws.Cells.Item(17, 31).Value = 32
$ws.Cells.Item(17, 32).Value = 20
$ws.Cells.Item(17, 33).Value = 13.5
$ws.Cells.Item(17, 34).Value = 970
$ws.Cells.Item(17, 35).Value = 42
$ws.Cells.Item(17, 36).Value = 44
$ws.Cells.Item(17, 37).Value = 32
$ws.Cells.Item(17, 38).Value = 42
$ws.Cells.Item(17, 39).Value = 100
$ws.Cells.Item(17, 40).Value = 29
$ws.Cells.Item(17, 41).Value = 25
# Row 18
$ws.Cells.Item(18, 6).Value = 1.38
$ws.Cells.Item(18, 7).Value = 1.46
$ws.Cells.Item(18, 8).Value = 8.4
$ws.Cells.Item(18, 9).Value = 10.5
$ws.Cells.Item(18, 10).Value = 5
$ws.Cells.Item(18, 11).Value = 5.9
$ws.Cells.Item(18, 13).Value = 1.03
$ws.Cells.Item(18, 14).Value = 4.8
$ws.Cells.Item(18, 15).Value = 1.2
$ws.Cells.Item(18, 18).Value = 1.55
$ws.Cells.Item(18, 19).Value = 2.56
$ws.Cells.Item(18, 20).Value = 1.88
$ws.Cells.Item(18, 21).Value = 1.94
$ws.Cells.Item(18, 22).Value = 1.11
$ws.Cells.Item(18, 23).Value = 3.15
$ws.Cells.Item(18, 24).Value = 24
$ws.Cells.Item(18, 25).Value = 34
$ws.Cells.Item(18, 26).Value = 85
$ws.Cells.Item(18, 27).Value = 330
$ws.Cells.Item(18, 28).Value = 10.5
$ws.Cells.Item(18, 29).Value = 13
$ws.Cells.Item(18, 30).Value = 34
$ws.Cells.Item(18, 31).Value = 150
$ws.Cells.Item(18, 32).Value = 9.800000000000001
$ws.Cells.Item(18, 33).Value = 11
$ws.Cells.Item(18, 34).Value = 26
$ws.Cells.Item(18, 35).Value = 120
$ws.Cells.Item(18, 36).Value = 13
$ws.Cells.Item(18, 37).Value = 15.5
$ws.Cells.Item(18, 38).Value = 34
$ws.Cells.Item(18, 39).Value = 150
$ws.Cells.Item(18, 40).Value = 5.8
$ws.Cells.Item(18, 41).Value = 170
